# Weekly data refresh: insert two new daily price rows (Arica - Zanahoria)
# right after the existing row for 2021-01-28, shifting every row below it
# down by two (old row 106 -> new row 108, ... old row 175 -> new row 177).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 106-107; everything from the old row 106 onward
# (previously through row 175) shifts down to rows 108-177.
$ws.Range("A106:A107").EntireRow.Insert()

# --- New row 106 ---
$ws.Range("A106").Value = 1
$ws.Range("B106").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C106").Value = "Arica y Parinacota"
$ws.Range("D106").Value = 44438
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = 15
$ws.Range("F106").Value = 100114013
$ws.Range("G106").Value = "Zanahoria"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 70
$ws.Range("K106").Value = 8000
$ws.Range("L106").Value = 9000
$ws.Range("M106").Value = 8500
$ws.Range("N106").Value = "$/saco 25 kilos"
$ws.Range("O106").Value = "Región de Arica y Parinacota"
$ws.Range("P106").Value = 340
$ws.Range("Q106").Value = 25
$ws.Range("R106").Value = "Hortaliza"

# --- New row 107 ---
$ws.Range("A107").Value = 1
$ws.Range("B107").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C107").Value = "Arica y Parinacota"
$ws.Range("D107").Value = 44438
$ws.Range("D107").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E107").Value = 15
$ws.Range("F107").Value = 100114013
$ws.Range("G107").Value = "Zanahoria"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Segunda"
$ws.Range("J107").Value = 50
$ws.Range("K107").Value = 8000
$ws.Range("L107").Value = 8500
$ws.Range("M107").Value = 8250
$ws.Range("N107").Value = "$/saco 25 kilos"
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 330
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = "Hortaliza"
